# Insert a new data row before the existing row 205 ("Vega Modelo de Temuco" /
# "Granada" price log). Excel shifts rows 205:301 down to 206:302 and keeps
# all of their values intact; we only need to populate the newly inserted
# row 205 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(205).Insert()

$ws.Range("A205").Value = 10
$ws.Range("B205").Value = "Vega Modelo de Temuco"
$ws.Range("C205").Value = "La Araucanía"
$ws.Range("D205").Value = 45202
$ws.Range("E205").Value = 9
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100104
$ws.Range("H205").Value = "Frutos de pepita"
$ws.Range("I205").Value = 100104001
$ws.Range("J205").Value = "Granada"
$ws.Range("K205").Value = "Wonderfull"
$ws.Range("L205").Value = "Primera"
$ws.Range("M205").Value = 80
$ws.Range("N205").Value = 17000
$ws.Range("O205").Value = 17000
$ws.Range("P205").Value = 17000
$ws.Range("Q205").Value = "$/bandeja 10 kilos granel"
$ws.Range("R205").Value = "Provincia de Limarí"
$ws.Range("S205").Value = 1700
$ws.Range("T205").Value = 10
